# GDP Growth Rates.xlsx -- apply commit "updated rmi files and update to 3.4.3"
#
# Summary of the real content/structural changes this script reproduces
# (cosmetic-only bits like internal style-table index renumbering and the
# raw pixel window geometry are left to the host to derive naturally):
#   - About: row 1 gets a new "Oregon" label in B1, and the date in C1 is
#     bumped from 2021-07-07 to 2022-09-20.
#   - Data: the stray formatting that had been painted across columns
#     E:J on several rows is cleared out (columns E/F only really carry
#     data through 2022 in row 3); row 11 (pure formatting, no data)
#     disappears once cleared.
#   - GDPGR-alternate: now references the new GDPGR-bau sheet instead of
#     reading the Data sheet directly, and becomes the active/selected
#     tab (selection parked on B3).
#   - GDPGR-bau / Data keep their formulas & values; only bookkeeping
#     (which sheet is active) changes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# About sheet
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Range("B1").Value = "Oregon"
$about.Range("C1").Value = 44824

# ---------------------------------------------------------------------
# Data sheet -- strip the leftover E:J formatting that doesn't belong
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")

# Row 3: E3 keeps its value (20767) but loses its formatting; F3:J3 were
# format-only and simply disappear once cleared.
$data.Range("E3").ClearFormats()
$data.Range("F3:J3").Clear()

# Rows 4, 5, 6, 8, 9, 11, 13: format-only cells in E:J -- clear them out
# entirely (row 11 had nothing else in it, so it vanishes completely).
$data.Range("E4:J4").Clear()
$data.Range("E5:J5").Clear()
$data.Range("E6:J6").Clear()
$data.Range("E8:J8").Clear()
$data.Range("E9:J9").Clear()
$data.Range("E11:J11").Clear()
$data.Range("E13:J13").Clear()

# ---------------------------------------------------------------------
# GDPGR-alternate sheet -- now sources its BAU-relative growth rate from
# the GDPGR-bau sheet (which itself reads Data!B13) instead of reading
# Data!B14 directly.
# ---------------------------------------------------------------------
$alt = $wb.Worksheets.Item("GDPGR-alternate")
$alt.Range("B2").Formula = "='GDPGR-bau'!B2"

# Make GDPGR-alternate the active/selected tab, with B3 selected.
$alt.Select()
$alt.Range("B3").Select()
